$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename statement labels in column A (auto increment public id naming convention)
$ws.Range("A2").Value = "statement-01"
$ws.Range("A3").Value = "statement-01"
$ws.Range("A4").Value = "statement-01"
$ws.Range("A5").Value = "statement-02"
$ws.Range("A6").Value = "statement-02"

# Clear the "Public ID" column values (auto-incremented now), keep formatting
$ws.Range("B2:B6").ClearContents()
